$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Concepts")

$xlPasteFormats = -4122
$xlPasteValues  = -4163

# ---- Metadata sheet ----------------------------------------------------
# Date property: bump to the new build timestamp.
$ws1.Range("B8").Value = "2024-09-13T20:57:31+00:00"

# Count property: 2 -> 3. Enter it as a formula returning the text "3"
# and then paste-special just the value back onto itself; this keeps the
# cell's text data type (matching the rest of the sheet) without leaving
# any left-over number-format/quote-prefix style behind.
$ws1.Range("B22").Formula = "=""3"""
$ws1.Range("B22").Copy() | Out-Null
$ws1.Range("B22").PasteSpecial($xlPasteValues) | Out-Null

# ---- Concepts sheet -----------------------------------------------------
# Add the new "unknown" / "Unknown" concept as row 4, matching row 3's
# layout and formatting (Level/Code/Display/Definition columns).
$ws2.Range("A3:D3").Copy() | Out-Null
$ws2.Range("A4:D4").PasteSpecial($xlPasteFormats) | Out-Null

# Level column must stay text "1" (same as rows 2-3) rather than a number.
$ws2.Range("A4").Formula = "=""1"""
$ws2.Range("A4").Copy() | Out-Null
$ws2.Range("A4").PasteSpecial($xlPasteValues) | Out-Null

$ws2.Range("B4").Value = "unknown"
$ws2.Range("C4").Value = "Unknown"

$excel.CutCopyMode = $false
